# "Update-Delete of Sumfone Tagadad"
#
# The sheet "RedeConflitGerTo_CentlitybByNod" (4th sheet, sheet4.xml) has a
# per-person row table in A1:Q131. Row 79 corresponds to the person
# "Sumfone Tagadad" (shared-string index 182). This row is removed entirely
# (an Excel "Delete Entire Row" on row 79), so every row below it shifts up
# by one, the used range shrinks from A1:Q131 to A1:Q130, and the
# sharedStrings reference count drops by one (the string itself stays in
# the table since it is still referenced from another sheet).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("RedeConflitGerTo_CentlitybByNod")

# Make sure we're working on/looking at the right sheet, like a user would.
$ws.Activate()

# Delete the entire row for "Sumfone Tagadad" (row 79) - this shifts every
# subsequent row up by one and shrinks the sheet's used range by one row.
$ws.Rows.Item(79).Delete()

# Reflect the post-edit cursor/selection position left behind in the file
# (activeCell="I82", sqref="I82" in the saved view state).
$ws.Range("I82").Select()
